# aggiornamento fino a 6/03
# Appends three new daily data rows (245-247) to the end of the existing
# table on Sheet1, continuing the date series in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: date (serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44319, 4, 14, 128.972823583602),
    @(44320, 1, 14, 128.972823583602),
    @(44321, 0, 14, 128.972823583602)
)

$lastRow = 244
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]

    # Match the date-cell style used throughout column A (bordered, bold,
    # centered, custom date/time number format) by copying the format
    # from the last existing row.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
